# The workbook's "Obc" sheet had its AC2 cell read in as a number (23)
# instead of the text value it should have been, and the whole dataset
# needed a second, sibling tab ("Sebc") so every service/category reads
# from its own sheet instead of being crammed together.

$wb  = $excel.ActiveWorkbook
$obc = $wb.Worksheets.Item("Obc")

# --- Fix the mis-typed cell: store "23" as text (quote-prefixed), not a number ---
$obc.Range("AC2").Value = "'23"

# --- Add the new "Sebc" tab as a duplicate of "Obc", placed right after it ---
$obc.Copy($null, $obc) | Out-Null
$sebc = $wb.Worksheets.Item("Obc (2)")
$sebc.Name = "Sebc"

# Point the new tab's selection at AD4 (it opens un-selected/not the active tab)
$sebc.Activate() | Out-Null
$sebc.Range("AD4").Select() | Out-Null

# --- Restore "Obc" as the active/selected tab, with its own cursor at AB9 ---
$obc.Activate() | Out-Null
$obc.Range("AB9").Select() | Out-Null
